$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (A:T) following Dr Hou's advice:
# 4 sending clusters (ECs, FAPs, M2, sCs) x 3 target clusters (ECs, M2, sCs)
# all for Selplg -> Selp, each with Ligand/Receptor-expressing cells = 3.
$data = @(
    @("ECs","Selplg","Selp","ECs",3,1,1.045864666666667,3.137594,0.01599546767026659,0.01599546767026659,3,1,15.960008,47.880024,0.9899803616776065,0.9899803616776066,16.69200844691733,150.228076022256,0.01583519886941298,0.01583519886941298),
    @("ECs","Selplg","Selp","M2",3,1,1.045864666666667,3.137594,0.01599546767026659,0.01599546767026659,2,0.6666666666666666,0.06197766666666666,0.185933,0.003844401134548353,0.003844401134548354,0.06482025168911111,0.583382265202,[double]"6.149299405920438E-05",[double]"6.149299405920438E-05"),
    @("ECs","Selplg","Selp","sCs",3,1,1.045864666666667,3.137594,0.01599546767026659,0.01599546767026659,1,0.3333333333333333,0.09955433333333334,0.298663,0.006175237187845165,0.006175237187845166,0.1041203596468889,0.937083236822,[double]"9.877580679440529E-05",[double]"9.877580679440531E-05"),
    @("FAPs","Selplg","Selp","ECs",3,1,2.006702333333334,6.020107,0.03069053130839924,0.03069053130839923,3,1,15.960008,47.880024,0.9899803616776065,0.9899803616776066,32.02698529361867,288.242867642568,0.03038302328476698,0.03038302328476698),
    @("FAPs","Selplg","Selp","M2",3,1,2.006702333333334,6.020107,0.03069053130839924,0.03069053130839923,2,0.6666666666666666,0.06197766666666666,0.185933,0.003844401134548353,0.003844401134548354,0.1243707283145556,1.119336554831,0.0001179867133819018,0.0001179867133819018),
    @("FAPs","Selplg","Selp","sCs",3,1,2.006702333333334,6.020107,0.03069053130839924,0.03069053130839923,1,0.3333333333333333,0.09955433333333334,0.298663,0.006175237187845165,0.006175237187845166,0.1997759129934445,1.797983216941,0.0001895213102503533,0.0001895213102503533),
    @("M2","Selplg","Selp","ECs",3,1,59.82593033333333,179.477791,0.914978548362651,0.914978548362651,3,1,15.960008,47.880024,0.9899803616776065,0.9899803616776066,954.8223267274427,8593.400940546984,0.9058107942353086,0.9058107942353087),
    @("M2","Selplg","Selp","M2",3,1,59.82593033333333,179.477791,0.914978548362651,0.914978548362651,2,0.6666666666666666,0.06197766666666666,0.185933,0.003844401134548353,0.003844401134548354,3.707871568222555,33.370844114003,0.003517544569412781,0.003517544569412781),
    @("M2","Selplg","Selp","sCs",3,1,59.82593033333333,179.477791,0.914978548362651,0.914978548362651,1,0.3333333333333333,0.09955433333333334,0.298663,0.006175237187845165,0.006175237187845166,5.955930610381444,53.603375493433,0.005650209557929628,0.005650209557929629),
    @("sCs","Selplg","Selp","ECs",3,1,2.506566,7.519698,0.03833545265868316,0.03833545265868316,3,1,15.960008,47.880024,0.9899803616776065,0.9899803616776066,40.004813412528,360.043320712752,0.03795134528811791,0.03795134528811792),
    @("sCs","Selplg","Selp","M2",3,1,2.506566,7.519698,0.03833545265868316,0.03833545265868316,2,0.6666666666666666,0.06197766666666666,0.185933,0.003844401134548353,0.003844401134548354,0.155351112026,1.398160008234,0.0001473768576944662,0.0001473768576944662),
    @("sCs","Selplg","Selp","sCs",3,1,2.506566,7.519698,0.03833545265868316,0.03833545265868316,1,0.3333333333333333,0.09955433333333334,0.298663,0.006175237187845165,0.006175237187845166,0.249539507086,2.245855563774,0.000236730512870778,0.0002367305128707781)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
